# Updated cryptos list with freshly scraped Price (D) and Volume(1h) (E) values.
# D-column prices are plain numeric-looking strings (e.g. "0.998", "316.68") that
# must stay stored as literal text (matching the sheet's existing inline-string cells),
# not get auto-converted to the Number type by Excel's normal data-entry inference.
# Route them through a formula -> Copy -> PasteSpecial(values) round trip so the literal
# text lands in the cell without Excel's smart number-detection and without touching the
# cell's number format/style (a plain .Value assignment of "0.998" would otherwise be
# stored as the number 0.998 and strip trailing zeros from values like "74.00").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-LiteralText($cell, $text) {
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

$excel.CutCopyMode = 0

$updates = @(
    @{ Row = 2; D = '42.776.43'; E = '  +0.74%  ' }
    @{ Row = 3; D = '2.298.12'; E = '  -0.32%  ' }
    @{ Row = 4; D = '0.998'; E = '  -0.40%  ' }
    @{ Row = 5; D = '316.68'; E = '  -0.53%  ' }
    @{ Row = 6; D = '104.39'; E = '  +0.21%  ' }
    @{ Row = 7; D = $null; E = '  -0.64%  ' }
    @{ Row = 8; D = $null; E = '  -0.04%  ' }
    @{ Row = 9; D = '0.604'; E = '  -1.32%  ' }
    @{ Row = 10; D = '39.62'; E = '  -0.83%  ' }
    @{ Row = 11; D = '0.0909'; E = '  -0.32%  ' }
    @{ Row = 12; D = '8.53'; E = '  +1.63%  ' }
    @{ Row = 13; D = $null; E = '  +2.48%  ' }
    @{ Row = 14; D = $null; E = '  +3.43%  ' }
    @{ Row = 15; D = $null; E = '  +0.16%  ' }
    @{ Row = 16; D = '2.643.40'; E = '  -0.54%  ' }
    @{ Row = 17; D = '2.292.89'; E = '  -0.33%  ' }
    @{ Row = 18; D = '42.669.97'; E = '  +0.07%  ' }
    @{ Row = 19; D = '15.12'; E = '  +36.39%  ' }
    @{ Row = 20; D = '7.54'; E = '  +0.20%  ' }
    @{ Row = 21; D = $null; E = '  +0.06%  ' }
    @{ Row = 22; D = '74.00'; E = '  +1.12%  ' }
    @{ Row = 23; D = $null; E = '  -1.28%  ' }
    @{ Row = 24; D = '266.27'; E = '  -5.22%  ' }
    @{ Row = 25; D = $null; E = '  -2.32%  ' }
    @{ Row = 26; D = $null; E = '  +0.49%  ' }
    @{ Row = 27; D = $null; E = '  +0.54%  ' }
    @{ Row = 28; D = '2.28'; E = '  -4.77%  ' }
    @{ Row = 29; D = '6.86'; E = '  +16.07%  ' }
    @{ Row = 30; D = '22.60'; E = '  -1.35%  ' }
    @{ Row = 31; D = '37.39'; E = '  +3.65%  ' }
    @{ Row = 32; D = '166.33'; E = '  +0.66%  ' }
    @{ Row = 33; D = $null; E = '  +0.22%  ' }
    @{ Row = 34; D = $null; E = '  -4.05%  ' }
    @{ Row = 35; D = $null; E = '  +0.66%  ' }
    @{ Row = 36; D = $null; E = '  -2.96%  ' }
    @{ Row = 37; D = $null; E = '  -1.79%  ' }
    @{ Row = 38; D = $null; E = '  -6.07%  ' }
    @{ Row = 39; D = '3.74'; E = '  -0.09%  ' }
    @{ Row = 40; D = '2.70'; E = '  -2.95%  ' }
    @{ Row = 41; D = $null; E = '  +4.48%  ' }
    @{ Row = 42; D = '70.43'; E = '  +0.50%  ' }
    @{ Row = 43; D = $null; E = '  +0.89%  ' }
    @{ Row = 44; D = '95.59'; E = '  -2.53%  ' }
    @{ Row = 45; D = $null; E = '  -0.32%  ' }
    @{ Row = 46; D = '12.36'; E = '  +1.74%  ' }
    @{ Row = 47; D = '115.18'; E = '  +2.62%  ' }
    @{ Row = 48; D = '80.05'; E = '  +0.14%  ' }
    @{ Row = 49; D = '1.706.36'; E = '  +6.26%  ' }
    @{ Row = 50; D = '8.81'; E = '  -1.71%  ' }
    @{ Row = 51; D = $null; E = '  -3.77%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        Set-LiteralText $ws.Cells.Item($u.Row, 4) $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

$excel.CutCopyMode = 0
